# Updated test cases for Sprint 1 stories.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new test case row (row 3) for PrintApp.002.
$ws.Range("A3").Value = "PrintApp.002"
$ws.Range("B3").Value = "Verify that staff members can print the appointment calendars by day."
$ws.Range("C3").Value = "1. Log in as a staff member and access the print calendar page.`n2. Attempt to print the calendar."

# Update the expected results for the existing PrintApp.001 test case (D2)
# to mention the default set of instructions.
$ws.Range("D2").Value = "2.a. The applicant should be able to successfully print out their appointment.`n2.b. The print out should contain the following information: Location, Time, Date, Office Address, Office Number, User's Name, User's Number, User's Email, and a default set of instructions."

$ws.Range("D3").Value = "2.a. The calendar for the day should be printed successfully.`n2.b. The print out should contain the following information: Location, Office Address, Applicant's Name, Applicant's Phone Number, Applicant's Email, Appointment Date, Appointment Time, and a default set of instructions."

# Match wrap-text styling used by the other data rows.
$ws.Range("A3:D3").WrapText = $true

# Adjust row heights to match the new content.
$ws.Rows.Item(2).RowHeight = 105
$ws.Rows.Item(3).RowHeight = 120

# Update the active cell selection to reflect where the user would be after
# entering the new row of data.
$ws.Range("D4").Select()
